$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates -----------------------------------------------------
# F2 (Exp Doc Nbr) already stored as text (Text number format) -> plain value assign keeps it text.
$ws.Range("F2").Value = "9229592017"

# AM2 (Func Loc) already stored as text (Text number format) -> plain value assign keeps it text.
$ws.Range("AM2").Value = "1640840748"

# New Quote Start Date / Quote End Date entries for row 2.
$ws.Range("DD2").Value = 42413
$ws.Range("DE2").Value = 42503

# --- Row 3 updates -----------------------------------------------------
# F3 has a Number format applied (numFmtId 1) even though the stored content is
# text, so toggle the format to Text while assigning, then restore the original
# Number format so the cell keeps its original appearance/style.
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "3340245657"
$ws.Range("F3").NumberFormat = "0"

# AM3 already stored as text (Text number format) -> plain value assign keeps it text.
$ws.Range("AM3").Value = "1640840748"

# New Quote Start Date / Quote End Date entries for row 3.
$ws.Range("DD3").Value = 42414
$ws.Range("DE3").Value = 42504

# --- View state ----------------------------------------------------------
$ws.Range("DE2:DE3").Select()
